# Actualización automática 2025-09-08 09:55:08
# Update the "CUMPLIMIENTO MENSUAL" sheet: remove the stale group rows
# (GRANITO, LED, PANELES PU, PANELES PVC) and refresh the remaining
# figures with the latest budget/sales numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Remove obsolete group rows (bottom-up so row numbers of the rows
# still to be removed don't shift while we work).
$ws.Rows.Item(14).Delete()   # PANELES PVC
$ws.Rows.Item(13).Delete()   # PANELES PU
$ws.Rows.Item(9).Delete()    # LED
$ws.Rows.Item(5).Delete()    # GRANITO

# Refresh the figures for every remaining group with the latest values.
$ws.Range("B2").Value = "240X120 PORCELANATO"
$ws.Range("C2").Value = 1377.24089543035
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1377.24089543035
$ws.Range("F2").Value = 0

$ws.Range("B3").Value = "240X80 PORCELANATO"
$ws.Range("C3").Value = 2582.41380675037
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2582.41380675037
$ws.Range("F3").Value = 0

$ws.Range("B4").Value = "FREGADEROS DE COCINA"
$ws.Range("C4").Value = 558.15203605817
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 558.15203605817
$ws.Range("F4").Value = 0

$ws.Range("B5").Value = "GRIFERIAS"
$ws.Range("C5").Value = 150
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 150
$ws.Range("F5").Value = 0

$ws.Range("B6").Value = "INODOROS"
$ws.Range("C6").Value = 918.796443341442
$ws.Range("D6").Value = 807
$ws.Range("E6").Value = 111.796443341442
$ws.Range("F6").Value = 0.8783229472081268

$ws.Range("B7").Value = "LAVABOS"
$ws.Range("C7").Value = 665.033262215681
$ws.Range("D7").Value = 160.89
$ws.Range("E7").Value = 504.143262215681
$ws.Range("F7").Value = 0.241927748792542

$ws.Range("B8").Value = "NO RESURTIBLES"
$ws.Range("C8").Value = 516.121873547834
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 516.121873547834
$ws.Range("F8").Value = 0

$ws.Range("B9").Value = "OTROS"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0

$ws.Range("B10").Value = "PANELES DECORATIVOS"
$ws.Range("C10").Value = 388.107983534392
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 388.107983534392
$ws.Range("F10").Value = 0

$ws.Range("B11").Value = "PIEDRA SINTERIZADA"
$ws.Range("C11").Value = 5844.44916370549
$ws.Range("D11").Value = 810.25
$ws.Range("E11").Value = 5034.19916370549
$ws.Range("F11").Value = 0.1386358196135436

$ws.Range("B12").Value = "PORCELANATO"
$ws.Range("C12").Value = 17675.3486842162
$ws.Range("D12").Value = 725.76
$ws.Range("E12").Value = 16949.5886842162
$ws.Range("F12").Value = 0.04106057611457996

$ws.Range("B13").Value = "PUERTAS DE SEGURIDAD"
$ws.Range("C13").Value = 364.412605947529
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 364.412605947529
$ws.Range("F13").Value = 0

$ws.Range("B14").Value = "SAL SOLUBLE"
$ws.Range("C14").Value = 667.683148387554
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 667.683148387554
$ws.Range("F14").Value = 0

$ws.Range("B15").Value = "TOTAL"
$ws.Range("C15").Value = 31707.75990313501
$ws.Range("D15").Value = 2503.9
$ws.Range("E15").Value = 29203.85990313501
$ws.Range("F15").Value = 0.07896805096447176

# Column width tweaks (D and F got a touch narrower).
# Excel's ColumnWidth is offset ~0.83 from the stored XML width units.
$ws.Columns.Item(4).ColumnWidth = 11.17
$ws.Columns.Item(6).ColumnWidth = 24.17
